$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values (column D) are plain decimal numbers (e.g. "199.98").
# The source data stores them as literal text (matching the "thousand.thousand" style
# entries like "76.457.03" elsewhere in the column), so force those specific cells to
# Text format first to avoid Excel auto-converting them into numeric values.
$textCells = @("D5", "D6", "D8", "D9", "D13", "D15", "D17", "D19", "D20", "D21", "D22", "D24", "D27", "D28", "D30", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D42", "D44", "D45", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "76.457.03"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.029.39"
$ws.Range("E3").Value = "  +3.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "199.98"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "627.83"
$ws.Range("E6").Value = "  +4.48%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "0.205"
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("D10").Value = "3.028.54"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "5.10"
$ws.Range("E13").Value = "  +4.57%  "
$ws.Range("D14").Value = "3.587.93"
$ws.Range("E14").Value = "  +4.18%  "
$ws.Range("D15").Value = "29.13"
$ws.Range("E15").Value = "  +6.11%  "
$ws.Range("D16").Value = "76.366.31"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "0.0000190"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "3.041.99"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("D19").Value = "13.46"
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("D20").Value = "9.06"
$ws.Range("E20").Value = "  +3.48%  "
$ws.Range("D21").Value = "373.09"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "4.35"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "73.24"
$ws.Range("E24").Value = "  +2.89%  "
$ws.Range("D25").Value = "3.192.57"
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "4.36"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +7.39%  "
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "508.89"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "1.95"
$ws.Range("E34").Value = "  +6.79%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").Value = "164.06"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "193.67"
$ws.Range("E38").Value = "  +7.09%  "
$ws.Range("D39").Value = "20.00"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").Value = "0.382"
$ws.Range("E40").Value = "  +10.13%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "0.111"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "42.50"
$ws.Range("E45").Value = "  +5.92%  "
$ws.Range("E46").Value = "  +4.86%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "0.716"
$ws.Range("E48").Value = "  +8.67%  "
$ws.Range("E49").Value = "  +5.11%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  +3.69%  "
